# Generate Report for Handoff
# Adds a new file "7d65ac3a-205f-4b2c-83a2-ea551215c72e.md" as row 8 to the
# Overview, zh-cn and de-de sheets/tables, following the same pattern as the
# existing "564bd0cf-42e9-4340-a0e9-fb94fd5c91e4.md" row (status: Ready for
# handoff / True, dependency False, no reference tokens etc).

$wb = $excel.ActiveWorkbook

$fileName   = "7d65ac3a-205f-4b2c-83a2-ea551215c72e.md"
$pathName   = "e2e\7d65ac3a-205f-4b2c-83a2-ea551215c72e.md"
$commitSha  = "714ebb150d80af156cbfa7cc82ee75ab27f3a527"
$zhXlf      = "7d65ac3a-205f-4b2c-83a2-ea551215c72e." + $commitSha + ".zh-cn.xlf"
$deXlf      = "7d65ac3a-205f-4b2c-83a2-ea551215c72e." + $commitSha + ".de-de.xlf"
$hoDate     = "2016-10-27 09:15:08"
$hoZhDate   = "2016-10-27 09:14:55"
$hoDeDate   = $hoDate
$ghUrl      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/" + $commitSha + "/e2e/" + $fileName

# ---------------------------------------------------------------------
# Overview sheet (sheet1) - columns: File Name, Path And Name, Extension,
# Publish URL, zh-cn, de-de, Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add()

$wsOverview.Range("A8").Value = $fileName
$wsOverview.Range("B8").Value = $pathName
$wsOverview.Hyperlinks.Add($wsOverview.Range("B8"), $ghUrl, "", "", $pathName)
$wsOverview.Range("C8").Value = ".md"
$wsOverview.Range("D8").Value = ""
$wsOverview.Range("E8").Value = "Ready for handoff"
$wsOverview.Range("F8").Value = "Ready for handoff"
$wsOverview.Range("G8").Value = $hoDate
$wsOverview.Range("G8").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# zh-cn sheet (sheet2) - columns: Source File Name, File Extension, Status,
# Source Path, Priority, Content Duplicate, Latest Handoff File,
# Latest Handoff Datetime, Latest Target File, Latest Handback File,
# Latest Handback DateTime, Reference Tokens, To be localized,
# Dependency From, Has metadata, Error Detail
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add()

$wsZhCn.Range("A8").Value = $fileName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A8"), $ghUrl, "", "", $fileName)
$wsZhCn.Range("B8").Value = ".md"
$wsZhCn.Range("C8").Value = "Ready for handoff"
$wsZhCn.Range("D8").Value = "e2e"
$wsZhCn.Range("E8").Value = "ht"
$wsZhCn.Range("F8").Value = "False"
$wsZhCn.Range("G8").Value = $zhXlf
$wsZhCn.Range("H8").Value = $hoZhDate
$wsZhCn.Range("H8").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I8").Value = ""
$wsZhCn.Range("J8").Value = ""
$wsZhCn.Range("K8").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K8").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L8").Value = ""
$wsZhCn.Range("M8").Value = "True"
$wsZhCn.Range("N8").Value = ""
$wsZhCn.Range("O8").Value = "False"
$wsZhCn.Range("P8").Value = ""

# ---------------------------------------------------------------------
# de-de sheet (sheet3) - same columns as zh-cn
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add()

$wsDeDe.Range("A8").Value = $fileName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A8"), $ghUrl, "", "", $fileName)
$wsDeDe.Range("B8").Value = ".md"
$wsDeDe.Range("C8").Value = "Ready for handoff"
$wsDeDe.Range("D8").Value = "e2e"
$wsDeDe.Range("E8").Value = "ht"
$wsDeDe.Range("F8").Value = "False"
$wsDeDe.Range("G8").Value = $deXlf
$wsDeDe.Range("H8").Value = $hoDeDate
$wsDeDe.Range("H8").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I8").Value = ""
$wsDeDe.Range("J8").Value = ""
$wsDeDe.Range("K8").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K8").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L8").Value = ""
$wsDeDe.Range("M8").Value = "True"
$wsDeDe.Range("N8").Value = ""
$wsDeDe.Range("O8").Value = "False"
$wsDeDe.Range("P8").Value = ""

Write-Output "Handoff report row added to Overview, zh-cn and de-de sheets."
